$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulated run details (log write mode update) - new values per commit
# Row 2
$ws.Cells.Item(2, 3).Value = 0.8949596881866455
$ws.Cells.Item(2, 5).Value = 5440.353396191823
$ws.Cells.Item(2, 8).Value = 0.1767605951137659
$ws.Cells.Item(2, 9).Value = 0.163262506580145
$ws.Cells.Item(2, 10).Value = 0.1583030468112386
$ws.Cells.Item(2, 11).Value = 0.1566933580605231
$ws.Cells.Item(2, 12).Value = 0.1557858736227853
$ws.Cells.Item(2, 13).Value = 0.1553796858760203
$ws.Cells.Item(2, 14).Value = 0.1551745850355836
$ws.Cells.Item(2, 15).Value = 0.1550705040120784
$ws.Cells.Item(2, 16).Value = 0.1550176867762698
$ws.Cells.Item(2, 17).Value = 0.1549908839998894
$ws.Cells.Item(2, 18).Value = 0.1549772825909799
$ws.Cells.Item(2, 19).Value = 0.1549705436789115
$ws.Cells.Item(2, 20).Value = 0.1549670692140769
$ws.Cells.Item(2, 21).Value = 0.1549652541951334
$ws.Cells.Item(2, 22).Value = 0.1549643060509092
$ws.Cells.Item(2, 23).Value = 0.1545057389291681
$ws.Cells.Item(2, 24).Value = 0.1544237506637995
$ws.Cells.Item(2, 25).Value = 0.1540497738049088

# Row 3
$ws.Cells.Item(3, 3).Value = 1.057798147201538
$ws.Cells.Item(3, 5).Value = 5039.478824497334
$ws.Cells.Item(3, 7).Value = 0.1933384005523712
$ws.Cells.Item(3, 8).Value = 0.1799754386507599
$ws.Cells.Item(3, 9).Value = 0.1680540835694761
$ws.Cells.Item(3, 10).Value = 0.1595033687699511
$ws.Cells.Item(3, 11).Value = 0.1516138664109001
$ws.Cells.Item(3, 12).Value = 0.1516138664109001
$ws.Cells.Item(3, 13).Value = 0.1516138664109001
$ws.Cells.Item(3, 14).Value = 0.1489436918043434
$ws.Cells.Item(3, 15).Value = 0.1489436918043434
$ws.Cells.Item(3, 16).Value = 0.1481312594162714
$ws.Cells.Item(3, 17).Value = 0.1475553677965185
$ws.Cells.Item(3, 18).Value = 0.1475553677965185
$ws.Cells.Item(3, 19).Value = 0.147264816024982
$ws.Cells.Item(3, 20).Value = 0.1469399677367001
$ws.Cells.Item(3, 21).Value = 0.1469399677367001
$ws.Cells.Item(3, 22).Value = 0.1469368220452329
$ws.Cells.Item(3, 23).Value = 0.146235454668564
$ws.Cells.Item(3, 24).Value = 0.146235454668564
$ws.Cells.Item(3, 25).Value = 0.146235454668564

# Row 4
$ws.Cells.Item(4, 3).Value = 0.890625
$ws.Cells.Item(4, 4).Value = 5
$ws.Cells.Item(4, 5).Value = 4971.846252555533
$ws.Cells.Item(4, 7).Value = 0.1938177082958615
$ws.Cells.Item(4, 8).Value = 0.1864078920415088
$ws.Cells.Item(4, 9).Value = 0.1556450001230484
$ws.Cells.Item(4, 10).Value = 0.1556450001230484
$ws.Cells.Item(4, 11).Value = 0.1527350601427934
$ws.Cells.Item(4, 12).Value = 0.1436263642244136
$ws.Cells.Item(4, 13).Value = 0.1403560766463956
$ws.Cells.Item(4, 14).Value = 0.1403560766463956
$ws.Cells.Item(4, 15).Value = 0.1403560766463956
$ws.Cells.Item(4, 16).Value = 0.1402530361378027
$ws.Cells.Item(4, 17).Value = 0.1396972180259811
$ws.Cells.Item(4, 18).Value = 0.1388288705955758
$ws.Cells.Item(4, 19).Value = 0.1380801770980697
$ws.Cells.Item(4, 20).Value = 0.1376418114432195
$ws.Cells.Item(4, 21).Value = 0.1373926163311942
$ws.Cells.Item(4, 22).Value = 0.1373152751542007
$ws.Cells.Item(4, 23).Value = 0.1372894322480053
$ws.Cells.Item(4, 24).Value = 0.1371331417193906
$ws.Cells.Item(4, 25).Value = 0.1369170809465016

# Row 5
$ws.Cells.Item(5, 3).Value = 0.9218471050262451
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 5319.276856685684
$ws.Cells.Item(5, 8).Value = 0.1735916782631247
$ws.Cells.Item(5, 9).Value = 0.1614733528139802
$ws.Cells.Item(5, 10).Value = 0.1608321234485468
$ws.Cells.Item(5, 11).Value = 0.1597570122895474
$ws.Cells.Item(5, 12).Value = 0.1576041170885637
$ws.Cells.Item(5, 13).Value = 0.1572659713685058
$ws.Cells.Item(5, 14).Value = 0.1569722157368286
$ws.Cells.Item(5, 15).Value = 0.156324517469674
$ws.Cells.Item(5, 16).Value = 0.1559737367604662
$ws.Cells.Item(5, 17).Value = 0.1543925697545493
$ws.Cells.Item(5, 18).Value = 0.15224635055469
$ws.Cells.Item(5, 19).Value = 0.15224635055469
$ws.Cells.Item(5, 20).Value = 0.15224635055469
$ws.Cells.Item(5, 21).Value = 0.15224635055469
$ws.Cells.Item(5, 22).Value = 0.1518360453834512
$ws.Cells.Item(5, 23).Value = 0.1518360453834512
$ws.Cells.Item(5, 24).Value = 0.1518360453834512
$ws.Cells.Item(5, 25).Value = 0.1516896073428008

# Row 6
$ws.Cells.Item(6, 3).Value = 0.9375262260437012
$ws.Cells.Item(6, 5).Value = 5158.287553907163
$ws.Cells.Item(6, 7).Value = 0.1956649310266082
$ws.Cells.Item(6, 8).Value = 0.1745157418118972
$ws.Cells.Item(6, 9).Value = 0.1654085712640405
$ws.Cells.Item(6, 10).Value = 0.1609103478398919
$ws.Cells.Item(6, 11).Value = 0.1521786897514743
$ws.Cells.Item(6, 12).Value = 0.1521786897514743
$ws.Cells.Item(6, 13).Value = 0.1521786897514743
$ws.Cells.Item(6, 14).Value = 0.1502431894814248
$ws.Cells.Item(6, 15).Value = 0.1502431894814248
$ws.Cells.Item(6, 16).Value = 0.1492353598039636
$ws.Cells.Item(6, 17).Value = 0.1492353598039636
$ws.Cells.Item(6, 18).Value = 0.1492353598039636
$ws.Cells.Item(6, 19).Value = 0.1492353598039636
$ws.Cells.Item(6, 20).Value = 0.1490159579315052
$ws.Cells.Item(6, 21).Value = 0.1487764446454088
$ws.Cells.Item(6, 22).Value = 0.1486253009961733
$ws.Cells.Item(6, 23).Value = 0.1486253009961733
$ws.Cells.Item(6, 24).Value = 0.1485972791523557
$ws.Cells.Item(6, 25).Value = 0.1485514143061825

# Row 7
$ws.Cells.Item(7, 3).Value = 0.9374721050262451
$ws.Cells.Item(7, 4).Value = 6
$ws.Cells.Item(7, 5).Value = 5013.580305049602
$ws.Cells.Item(7, 7).Value = 0.1956649310266082
$ws.Cells.Item(7, 8).Value = 0.1803251139992895
$ws.Cells.Item(7, 9).Value = 0.1629978428180135
$ws.Cells.Item(7, 10).Value = 0.1555437479710646
$ws.Cells.Item(7, 11).Value = 0.1555437479710646
$ws.Cells.Item(7, 12).Value = 0.1540546227523222
$ws.Cells.Item(7, 13).Value = 0.1513452045070048
$ws.Cells.Item(7, 14).Value = 0.1488925573226196
$ws.Cells.Item(7, 15).Value = 0.1473057486562785
$ws.Cells.Item(7, 16).Value = 0.1473057486562785
$ws.Cells.Item(7, 17).Value = 0.1466799783034419
$ws.Cells.Item(7, 18).Value = 0.1466799783034419
$ws.Cells.Item(7, 19).Value = 0.1466799783034419
$ws.Cells.Item(7, 20).Value = 0.1459874589744318
$ws.Cells.Item(7, 21).Value = 0.1459874589744318
$ws.Cells.Item(7, 22).Value = 0.1459874589744318
$ws.Cells.Item(7, 23).Value = 0.1459393977084775
$ws.Cells.Item(7, 24).Value = 0.145900343439188
$ws.Cells.Item(7, 25).Value = 0.145730610234885

# Row 8
$ws.Cells.Item(8, 3).Value = 0.9062750339508057
$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(8, 5).Value = 5031.814944307063
$ws.Cells.Item(8, 7).Value = 0.1956649310266082
$ws.Cells.Item(8, 8).Value = 0.1645778924473497
$ws.Cells.Item(8, 9).Value = 0.1639963975265227
$ws.Cells.Item(8, 10).Value = 0.1553880447867191
$ws.Cells.Item(8, 11).Value = 0.1535724004106633
$ws.Cells.Item(8, 12).Value = 0.1508117815163711
$ws.Cells.Item(8, 13).Value = 0.1508117815163711
$ws.Cells.Item(8, 14).Value = 0.145455297840117
$ws.Cells.Item(8, 15).Value = 0.1444821195840719
$ws.Cells.Item(8, 16).Value = 0.1428593814968889
$ws.Cells.Item(8, 17).Value = 0.1399482020620552
$ws.Cells.Item(8, 18).Value = 0.1380860612925353
$ws.Cells.Item(8, 19).Value = 0.1380860612925353
$ws.Cells.Item(8, 20).Value = 0.1380860612925353
$ws.Cells.Item(8, 21).Value = 0.1380860612925353
$ws.Cells.Item(8, 22).Value = 0.1380860612925353
$ws.Cells.Item(8, 23).Value = 0.1380860612925353
$ws.Cells.Item(8, 24).Value = 0.1380860612925353
$ws.Cells.Item(8, 25).Value = 0.1380860612925353

# Row 9
$ws.Cells.Item(9, 3).Value = 0.9218747615814209
$ws.Cells.Item(9, 5).Value = 5070.0895054556
$ws.Cells.Item(9, 7).Value = 0.1956649310266082
$ws.Cells.Item(9, 8).Value = 0.1734577395176391
$ws.Cells.Item(9, 9).Value = 0.1663442636233587
$ws.Cells.Item(9, 10).Value = 0.1583424710489685
$ws.Cells.Item(9, 11).Value = 0.1554586994347532
$ws.Cells.Item(9, 12).Value = 0.1553274285010321
$ws.Cells.Item(9, 13).Value = 0.1516863431370425
$ws.Cells.Item(9, 14).Value = 0.1516863431370425
$ws.Cells.Item(9, 15).Value = 0.1516863431370425
$ws.Cells.Item(9, 16).Value = 0.1508586721296964
$ws.Cells.Item(9, 17).Value = 0.1500770146522692
$ws.Cells.Item(9, 18).Value = 0.1500354298476126
$ws.Cells.Item(9, 19).Value = 0.1488582536780025
$ws.Cells.Item(9, 20).Value = 0.1487915780229555
$ws.Cells.Item(9, 21).Value = 0.1481429111516316
$ws.Cells.Item(9, 22).Value = 0.1480904194804329
$ws.Cells.Item(9, 23).Value = 0.1473483447188515
$ws.Cells.Item(9, 24).Value = 0.1471836717231068
$ws.Cells.Item(9, 25).Value = 0.1468321541024483

# Row 10
$ws.Cells.Item(10, 3).Value = 0.7656238079071045
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 4623.705380305699
$ws.Cells.Item(10, 7).Value = 0.1929649790401971
$ws.Cells.Item(10, 8).Value = 0.1858876739886031
$ws.Cells.Item(10, 9).Value = 0.1723251597758633
$ws.Cells.Item(10, 10).Value = 0.1544909325568811
$ws.Cells.Item(10, 11).Value = 0.1461016346414996
$ws.Cells.Item(10, 12).Value = 0.1409746712075019
$ws.Cells.Item(10, 13).Value = 0.1355059157275773
$ws.Cells.Item(10, 14).Value = 0.1205822195364806
$ws.Cells.Item(10, 15).Value = 0.1141534872114345
$ws.Cells.Item(10, 16).Value = 0.1079924447369404
$ws.Cells.Item(10, 17).Value = 0.1033483520201393
$ws.Cells.Item(10, 18).Value = 0.1029538891833664
$ws.Cells.Item(10, 19).Value = 0.1004896582884087
$ws.Cells.Item(10, 20).Value = 0.1004361870568431
$ws.Cells.Item(10, 21).Value = 0.09985517860787904
$ws.Cells.Item(10, 22).Value = 0.09826746735873401
$ws.Cells.Item(10, 23).Value = 0.09826746735873401
$ws.Cells.Item(10, 24).Value = 0.09826746735873401
$ws.Cells.Item(10, 25).Value = 0.09813070916775241

# Row 11
$ws.Cells.Item(11, 3).Value = 0.953125
$ws.Cells.Item(11, 5).Value = 5132.022502926322
$ws.Cells.Item(11, 8).Value = 0.1864078920415088
$ws.Cells.Item(11, 9).Value = 0.1809165505885315
$ws.Cells.Item(11, 10).Value = 0.1666030053541784
$ws.Cells.Item(11, 11).Value = 0.160504440436756
$ws.Cells.Item(11, 12).Value = 0.1524174849973258
$ws.Cells.Item(11, 13).Value = 0.1506766678861566
$ws.Cells.Item(11, 14).Value = 0.1500024857470708
$ws.Cells.Item(11, 15).Value = 0.1484065752447079
$ws.Cells.Item(11, 16).Value = 0.1484065752447079
$ws.Cells.Item(11, 17).Value = 0.1484065752447079
$ws.Cells.Item(11, 18).Value = 0.1484065752447079
$ws.Cells.Item(11, 19).Value = 0.1484065752447079
$ws.Cells.Item(11, 20).Value = 0.1484065752447079
$ws.Cells.Item(11, 21).Value = 0.1484065752447079
$ws.Cells.Item(11, 22).Value = 0.1484065752447079
$ws.Cells.Item(11, 23).Value = 0.1484065752447079
$ws.Cells.Item(11, 24).Value = 0.1480394250083104
$ws.Cells.Item(11, 25).Value = 0.1480394250083104

